# Commit - Desktop Lab, 05/07/2024 as 18:41h.
# Reorders the "Resíduo" columns (Min X / Max X now immediately follow DP X,
# before DP Z / Resíduo Min Z / Resíduo Max Z) and refreshes the experimental
# data with newly measured values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (row 1) in their new order --------------------------------
$ws.Range("A1").Value = "Distância (m)"
$ws.Range("B1").Value = "DP X (mm)"
$ws.Range("C1").Value = "Resíduo Min X (mm)"
$ws.Range("D1").Value = "Resíduo Max X (mm)"
$ws.Range("E1").Value = "DP Z (mm)"
$ws.Range("F1").Value = "Resíduo Min Z (mm)"
$ws.Range("G1").Value = "Resíduo Max Z (mm)"

# --- Data rows (2-7) with the refreshed experimental values -------------
$data = @(
    @(1, 0.080000000000000002, -0.20999999999999999, 0.14999999999999999, 0.33000000000000002, -1.0700000000000001, 1.0800000000000001),
    @(2, 0.14999999999999999, -0.40000000000000002, 0.31, 0.88, -2.23, 2.3900000000000001),
    @(3, 0.23000000000000001, -0.58999999999999997, 0.46000000000000002, 1.1699999999999999, -3.8799999999999999, 3.3300000000000001),
    @(4, 0.28999999999999998, -0.76000000000000001, 0.63, 3.1699999999999999, -7.3799999999999999, 8.4299999999999997),
    @(5, 0.35999999999999999, -1.03, 0.79000000000000004, 3.2599999999999998, -9.9900000000000002, 12.640000000000001),
    @(6, 0.42999999999999999, -1.1299999999999999, 0.95999999999999996, 5.54, -17.07, 16.75)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $r++
}

# --- Column widths follow the header that now lives in each column ------
# (ColumnWidth is specified in "before padding" character units; Excel adds
#  ~5/6 of a character of internal padding when it stores the column width,
#  so the values below are chosen so the stored/serialized width lands on
#  the target width used by the other columns' headers.)
$ws.Columns.Item(3).ColumnWidth = 18.022135416666668
$ws.Columns.Item(4).ColumnWidth = 18.307291666666668
$ws.Columns.Item(5).ColumnWidth = 9.451822916666666
$ws.Columns.Item(6).ColumnWidth = 17.877604166666668
